$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview!G3 - "Latest HO Xliff Generate Date" for 438d8a02...
$wsOverview.Range("G3").Value = "2016-09-06 04:22:36"

# zh-cn sheet, row 3 (438d8a02 entry)
$wsZhCn.Range("H3").Value = "2016-09-06 04:22:25"
$wsZhCn.Range("K3").Value = "2016-09-06 04:23:15"

# de-de sheet, row 3 (438d8a02 entry)
$wsDeDe.Range("H3").Value = "2016-09-06 04:22:36"
$wsDeDe.Range("K3").Value = "2016-09-06 04:23:33"
